$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 100: Turno A, Operador 8524, manutenção ainda aguardando (reuses existing note)
$ws.Range("A100").Value = 46046
$ws.Range("B100").Value = "A"
$ws.Range("D100").Value = 8524
$ws.Range("F100").Value = " "
$ws.Range("G100").Value = "NÃO HOUVE PRODUÇÃO, AGUARDANDO MANUTENÇÃO DA MAQUINA"
$ws.Range("H100").Value = " "

# Row 101: Turno A, Operador 8524, manutenção concluída (novo texto)
$ws.Range("A101").Value = 46048
$ws.Range("B101").Value = "A"
$ws.Range("D101").Value = 8524
$ws.Range("F101").Value = " "
$ws.Range("G101").Value = "MANUTENÇÃO CONCLUÍDA POR VOLTA DAS 10:30, OPERADOR FICOU NAS MAQUINAS E PASSANDO PRIMER POIS O HORÁRIO JÁ ESTAVA PRÓXIMO AO DO ALMOÇO"
$ws.Range("H101").Value = " "

# Row 102: Turno B, OS 1406169, produção normal
$ws.Range("A102").Value = 46048
$ws.Range("B102").Value = "B"
$ws.Range("C102").Value = 1406169
$ws.Range("D102").Value = 8502
$ws.Range("E102").Value = 12097
$ws.Range("F102").Value = 10

# Match the saved file's view state: scrolled down near the new rows with
# F102 as the active cell.
$ws.Range("F102").Select()

$wb.Save()
